# Reorders the news-story rows (title/timestamp/historical distance/uri travel
# together as a unit) to reflect the newly added JSON record used for the
# time-bucket analysis, per the commit "added one json for time bucket analysis".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row content (columns A-E). Column D ("day_31_beyond") is unchanged for
# every row, so it is included for completeness but stays constant.
$rows = @(
    @{ Row = 2; A = "Storm Finally Exits a Snow-Covered East Coast";         B = "2009-12-20T15:24:31UTC"; C = 353; D = "day_31_beyond"; E = "https://www.nytimes.com/2009/12/21/us/21storm.html/?hp" },
    @{ Row = 3; A = "Record-breaking storm closes US federal government";    B = "2009-12-21T00:00:00UTC"; C = 354; D = "day_31_beyond"; E = "https://web.archive.org/web/20091222153852/http://news.yahoo.com/s/afp/20091221/ts_alt_afp/usweatherstorm" },
    @{ Row = 4; A = "Shoppers catch a break";                                B = "2009-12-21T13:43:10UTC"; C = 354; D = "day_31_beyond"; E = "https://lfpress.com/news/london/2009/12/19/12214646.html" },
    @{ Row = 5; A = "VDOT Continues to Address Winter Storm";                B = "2009-12-19T00:00:00UTC"; C = 352; D = "day_31_beyond"; E = "http://www.nbc29.com/Global/story.asp?S=11704169" },
    @{ Row = 6; A = "Crippling U.S. storm moves north";                     B = "2009-12-20T00:00:00UTC"; C = 353; D = "day_31_beyond"; E = "https://www.cbc.ca/news/world/crippling-u-s-storm-moves-north-1.787167" },
    @{ Row = 7; A = "Five dead as snowstorm engulfs US East Coast";          B = "2009-12-20T06:20:02UTC"; C = 353; D = "day_31_beyond"; E = "http://news.bbc.co.uk/2/hi/americas/8422652.stm" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
}

# The "uri" column cells carry live hyperlinks; rebuild them so each link
# target follows its row's new content. Per-item Hyperlinks.Delete() is a
# no-op in this host, but clearing the whole collection works, so drop every
# hyperlink and re-add the six with their correct (possibly new) targets.
$ws.Hyperlinks.Delete()
foreach ($r in $rows) {
    $ws.Hyperlinks.Add($ws.Range("E$($r.Row)"), $r.E)
}
